$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row ---
$ws.Range("B1").Value = "Total"
$ws.Range("C1").Value = "Evolução Total (%)"
$ws.Range("C1").Copy() | Out-Null
$ws.Range("D1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("D1").Value = "Qtd Produtos"
$excel.CutCopyMode = 0

# Row 2: 2020 (brand-new row, no Evolucao value)
$ws.Range("A2").Value = 2020
$ws.Range("B2").Value = 23866.41
$ws.Range("D2").Value = 262

# Row 3: 2021
$ws.Range("A3").Value = 2021
$ws.Range("B3").Value = 279918.14
$ws.Range("C3").Value = 1072.853981809581
$ws.Range("D3").Value = 2168

# Row 4: 2022
$ws.Range("A4").Value = 2022
$ws.Range("B4").Value = 504349.28
$ws.Range("C4").Value = 80.17741901257274
$ws.Range("D4").Value = 2708

# Row 5: 2023
$ws.Range("A5").Value = 2023
$ws.Range("B5").Value = 802464.89
$ws.Range("C5").Value = 59.1089591721039
$ws.Range("D5").Value = 3337

# Row 6: 2024
$ws.Range("A6").Value = 2024
$ws.Range("B6").Value = 1197062.47
$ws.Range("C6").Value = 49.17318937156241
$ws.Range("D6").Value = 4742

# Row 7: 2025
$ws.Range("A7").Value = 2025
$ws.Range("B7").Value = 396607.32
$ws.Range("C7").Value = -66.86828549557652
$ws.Range("D7").Value = 1635
